# Update the "department" column (C2:C7) on the "courses" sheet.
# Previously every row shared the single value "FACULTY OF HOSPITALITY".
# It is now split into two more specific department names:
#   - rows 2-4 (Certificate III / IV / Diploma fast-track)           -> "Hospitality"
#   - rows 5-7 (the bundled package rows)                            -> "Packages"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value = "Hospitality"
$ws.Range("C5:C7").Value = "Packages"

# Match the author's final selection/cursor position recorded in the file.
$ws.Range("D10").Select()
